# translating subs now uses sub length if available
#
# Adds a "model" column (B) recording "gpt-4o" for every video, shifting
# the previous B:K / M:N data right into C:L / N:O, and records three more
# translated videos (gakumas lilja hakusen / kamurogiku / sakura
# photograph) as new rows 7-9.
#
# Values are written directly into their final positions (rather than via
# a genuine column-insert) so the sheet's two pre-existing <col> overrides
# - the "#,##0" style pinned to physical column F and the "0.00" style
# pinned to physical column G - stay put and simply pick up new content,
# exactly like the authored workbook. Writes are sequenced so brand-new
# strings are first used in the same order the author typed them (new
# row 7 before the "model" column, then rows 8-9) so the shared-string
# table comes out in the expected order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Step 1: new video, row 7, typed in the *old* layout (A.."K", no model
# column yet) - this is what introduces "gakumas lilja hakusen" first.
# ---------------------------------------------------------------------
$ws.Range("A7").Value = "gakumas lilja hakusen"
$ws.Range("K7").Value = 0           # placeholder in old "cost_per_line" col, overwritten below

# ---------------------------------------------------------------------
# Step 2: shift the header row and data rows 1-7 one column to the right,
# B:K -> C:L (values/formulas only - done in final positions directly).
# ---------------------------------------------------------------------

# Header row
$ws.Range("C1").Value = "hours"
$ws.Range("D1").Value = "mins"
$ws.Range("E1").Value = "sec"
$ws.Range("F1").Value = "lines"
$ws.Range("G1").Value = "character_count"
$ws.Range("G1").NumberFormat = "#,##0"
$ws.Range("H1").Value = "average_character_count"
$ws.Range("I1").Value = "context_size"
$ws.Range("J1").Value = "cost_usd"
$ws.Range("K1").Value = "cost_per_character"
$ws.Range("L1").Value = "cost_per_line"

# Row 2 - gakumas lilja bond step 1
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 23
$ws.Range("E2").Value = 36
$ws.Range("F2").Value = 935
$ws.Range("G2").Value = 15967
$ws.Range("G2").NumberFormat = "#,##0"
$ws.Range("H2").Formula = "=G2/F2"
$ws.Range("H2").NumberFormat = "0.00"
$ws.Range("I2").Value = 3
$ws.Range("J2").Formula = "=4.98-2.58"
$ws.Range("J2").NumberFormat = "0.00"
$ws.Range("K2").Formula = "=J2/G2"
$ws.Range("K2").NumberFormat = "0.000000"
$ws.Range("L2").Formula = "=J2/F2"

# Row 3 - gakumas lilja bond step 2
$ws.Range("C3").Value = 0
$ws.Range("F3").Value = 705
$ws.Range("G3").Value = 10215
$ws.Range("G3").NumberFormat = "#,##0"
$ws.Range("H3").Formula = "=G3/F3"
$ws.Range("H3").NumberFormat = "0.00"
$ws.Range("I3").Value = 3
$ws.Range("J3").Value = 1.55
$ws.Range("J3").NumberFormat = "0.00"
$ws.Range("K3").Formula = "=J3/G3"
$ws.Range("K3").NumberFormat = "0.000000"
$ws.Range("L3").Formula = "=J3/F3"

# Row 4 - makeine vol2.1
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 10
$ws.Range("E4").Value = 45
$ws.Range("F4").Value = 331
$ws.Range("G4").Value = 4569
$ws.Range("G4").NumberFormat = "#,##0"
$ws.Range("H4").Formula = "=G4/F4"
$ws.Range("H4").NumberFormat = "0.00"
$ws.Range("I4").Value = 3
$ws.Range("J4").Value = 0.92
$ws.Range("K4").Formula = "=J4/G4"
$ws.Range("K4").NumberFormat = "0.000000"
$ws.Range("L4").Formula = "=J4/F4"

# Row 5 - makeine vol2.2
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 8
$ws.Range("E5").Value = 35
$ws.Range("F5").Value = 248
$ws.Range("G5").Value = 3496
$ws.Range("G5").NumberFormat = "#,##0"
$ws.Range("H5").Formula = "=G5/F5"
$ws.Range("H5").NumberFormat = "0.00"
$ws.Range("I5").Value = 3
$ws.Range("J5").Value = 0.77
$ws.Range("K5").Formula = "=J5/G5"
$ws.Range("K5").NumberFormat = "0.000000"
$ws.Range("L5").Formula = "=J5/F5"

# Row 6 - makeine vol2.3
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 11
$ws.Range("E6").Value = 26
$ws.Range("F6").Value = 351
$ws.Range("G6").Value = 4351
$ws.Range("G6").NumberFormat = "#,##0"
$ws.Range("H6").Formula = "=G6/F6"
$ws.Range("H6").NumberFormat = "0.00"
$ws.Range("I6").Value = 3
$ws.Range("J6").Value = 0.98
$ws.Range("K6").Formula = "=J6/G6"
$ws.Range("K6").NumberFormat = "0.000000"
$ws.Range("L6").Formula = "=J6/F6"

# Row 7 - gakumas lilja hakusen (finish shifting into C:L)
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 3
$ws.Range("E7").Value = 57
$ws.Range("F7").Value = 166
$ws.Range("G7").Value = 2283
$ws.Range("G7").NumberFormat = "#,##0"
$ws.Range("H7").Formula = "=G7/F7"
$ws.Range("H7").NumberFormat = "0.00"
$ws.Range("I7").Value = 3
$ws.Range("J7").Value = 0.38
$ws.Range("K7").Formula = "=J7/G7"
$ws.Range("K7").NumberFormat = "0.000000"
$ws.Range("L7").Formula = "=J7/F7"

# Helper table, old M2:N3 -> new N2:O3
$ws.Range("N2").Value = "avg_line_cost:"
$ws.Range("O2").Formula = "=AVERAGE(L:L)"
$ws.Range("N3").Value = "total lines:"
$ws.Range("O3").Formula = "=SUM(F:F)"
$ws.Range("M2").Value = ""
$ws.Range("M3").Value = ""

# ---------------------------------------------------------------------
# Step 3: new "model" column B - this is what introduces "model" then
# "gpt-4o" into the shared-string table (rows 1-7 first).
# ---------------------------------------------------------------------
$ws.Range("B1").Value = "model"
for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 2).Value = "gpt-4o"
}

# ---------------------------------------------------------------------
# Step 4: two brand new rows, 8 and 9, already in the final layout -
# these introduce "gakumas lilja kamurogiku" and "... sakura photograph".
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "gakumas lilja kamurogiku"
$ws.Range("B8").Value = "gpt-4o"
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 4
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 162
$ws.Range("F8").NumberFormat = "#,##0"
$ws.Range("G8").Value = 2667
$ws.Range("G8").NumberFormat = "#,##0"
$ws.Range("H8").Formula = "=G8/F8"
$ws.Range("H8").NumberFormat = "0.00"
$ws.Range("I8").Value = 3
$ws.Range("J8").Value = 0.38
$ws.Range("K8").Formula = "=J8/G8"
$ws.Range("K8").NumberFormat = "0.000000"
$ws.Range("L8").Formula = "=J8/F8"

$ws.Range("A9").Value = "gakumas lilja sakura photograph"
$ws.Range("B9").Value = "gpt-4o"
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 9
$ws.Range("E9").Value = 42
$ws.Range("F9").Value = 260
$ws.Range("F9").NumberFormat = "#,##0"
$ws.Range("G9").Value = 4231
$ws.Range("G9").NumberFormat = "#,##0"
$ws.Range("H9").Formula = "=G9/F9"
$ws.Range("H9").NumberFormat = "0.00"
$ws.Range("I9").Value = 3
$ws.Range("J9").Value = 0.81
$ws.Range("K9").Formula = "=J9/G9"
$ws.Range("K9").NumberFormat = "0.000000"
$ws.Range("L9").Formula = "=J9/F9"

# ---------------------------------------------------------------------
# Cosmetics: column B auto-fit width, final selection.
# ---------------------------------------------------------------------
$ws.Columns("B:B").AutoFit()
$ws.Range("N9").Select()
